# Generate Report for Handback
#
# For the "046f2de5-455d-4413-96fd-3b067f12f358" row (row 6) on both the
# zh-cn and de-de status sheets, the localization engine has produced a
# handback report: the "Latest Target File" (I), "Latest Handback File" (J)
# and "Latest Handback DateTime" (K) columns get populated, and the
# "Error Detail" (P) column records that the handback is stale. The
# "Error Detail" column is also widened so the message is readable.

$wb = $excel.ActiveWorkbook

$latestTargetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9dd1ee0d5f48107f0d7a708b8367f86eb7a98490/e2e/046f2de5-455d-4413-96fd-3b067f12f358.md"
$latestTargetDisplay = "046f2de5-455d-4413-96fd-3b067f12f358.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8f48bd55f51b3e838ee827a6f0288eb0d1e90449/e2e/046f2de5-455d-4413-96fd-3b067f12f358.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9dd1ee0d5f48107f0d7a708b8367f86eb7a98490/e2e/046f2de5-455d-4413-96fd-3b067f12f358.md."

function Update-StatusSheet($SheetName, $HandoffXlf, $HandbackDateTime) {
    $ws = $wb.Worksheets.Item($SheetName)

    # Widen the "Error Detail" column (P / column 16) so the new message is visible.
    $ws.Columns.Item(16).ColumnWidth = 39.17

    # Latest Target File: a hyperlink to the latest handed-off markdown source.
    $ws.Range("I6").Value = $latestTargetDisplay
    $ws.Hyperlinks.Add($ws.Range("I6"), $latestTargetUrl, "", "", $latestTargetDisplay)

    # Latest Handback File for this language.
    $ws.Range("J6").Value = $HandoffXlf

    # Latest Handback DateTime for this language.
    $ws.Range("K6").Value = $HandbackDateTime

    # Error Detail: handback file version mismatch message.
    $ws.Range("P6").Value = $errorDetail
}

Update-StatusSheet "zh-cn" "046f2de5-455d-4413-96fd-3b067f12f358.681076b0591596dd3990540d92be010845cf451e.zh-cn.xlf" "2016-08-28 10:51:47"
Update-StatusSheet "de-de" "046f2de5-455d-4413-96fd-3b067f12f358.681076b0591596dd3990540d92be010845cf451e.de-de.xlf" "2016-08-28 10:51:54"
